# Append new scraped listings (2025-12-24 12:38 JST run) to the "ランサーズ" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove existing hyperlinks first so we can rebuild them cleanly in the
# correct order once the new rows are in place.
$ws.Hyperlinks.Delete()

# Widen the price column slightly (30 -> 32), as in the target workbook.
# Note: the ColumnWidth property is offset from the raw stored column width
# by ~0.8333 (5/6) character units in this engine (matches classic Excel's
# character-width/MDW quirk), so subtract that to land on exactly 32.
$ws.Columns.Item(4).ColumnWidth = 31.16666667

$timestamp = "2025-12-24 12:38:44"

# Full refreshed dataset (11 listings), newest run, same timestamp for all
# rows, as produced by the scraper for this run.
$titles = @(
  "産業機械向けAI異常検知・状態推定システムの開発・導入支援エンジニア募集(AI/エッジ・組み込み)",
  "【フルタイム】最先端AI(LLM)開発エンジニア募集!新規プロダクトの核となる開発パートナーを募集",
  "【週5日】法人向け生成AIサービス(RAG・議事録機能)のコア開発を担うリードエンジニア募集",
  "【急募】AI活用でPDFタイトル修正のフリーランス募集!",
  "【急募】野球スコアボードシステム開発のフリーランス募集",
  "初回 【AWSクラウドリフト】業務アプリ移行支援エンジニア募集(Java / .NET)",
  "現品票管理・納品書・請求書のシステムづくり",
  "急募 限定公開 限定公開の仕事",
  "【電卓設計】ハードウェアとソフトウェアの専門家を募集!",
  "【電卓設計】ハードウェアとソフトウェアの専門家を募集!",
  "【急募】お名前VPSでのFTP・WPファイルアップロード改善依頼"
)

$categories = @(
  "システム開発",
  "システム開発",
  "システム開発",
  "システム開発",
  "システム開発",
  "システム開発",
  "システム開発",
  "システム開発",
  "システム開発",
  "システム開発",
  "システム開発"
)

$prices = @(
  "200,000 円 ~ 300,000 円 / 固定",
  "1,000,000 円 ~ 3,000,000 円 / 固定",
  "500,000 円 ~ 1,000,000 円 / 固定",
  "500,000 円 ~ 1,000,000 円 / 固定",
  "200,000 円 ~ 300,000 円 / 固定",
  "200,000 円 ~ 300,000 円 / 固定",
  "200,000 円 ~ 300,000 円 / 固定",
  "200,000 円 ~ 300,000 円 / 固定",
  "50,000 円 ~ 100,000 円 / 固定",
  "50,000 円 ~ 100,000 円 / 固定",
  "5,000 円 ~ 10,000 円 / 固定"
)

$deadlines = @(
  "期限情報なし",
  "期限情報なし",
  "期限情報なし",
  "期限情報なし",
  "期限情報なし",
  "期限情報なし",
  "期限情報なし",
  "期限情報なし",
  "期限情報なし",
  "期限情報なし",
  "期限情報なし"
)

$urls = @(
  "https://www.lancers.jp/work/detail/5450864",
  "https://www.lancers.jp/work/detail/5460294",
  "https://www.lancers.jp/work/detail/5460267",
  "https://www.lancers.jp/work/detail/5459721",
  "https://www.lancers.jp/work/detail/5459984",
  "https://www.lancers.jp/work/detail/5459847",
  "https://www.lancers.jp/work/detail/5459942",
  "https://www.lancers.jp/work/detail/5460299",
  "https://www.lancers.jp/work/detail/5459773",
  "https://www.lancers.jp/work/detail/5459232",
  "https://www.lancers.jp/work/detail/5459964"
)

$scores = @(383, 375, 375, 310, 118, 103, 53, 18, 18, 18, 10)

$skills = @(
  "🔥AI,Ai ◆開発",
  "🔥AI,Ai ◆開発",
  "🔥AI,Ai ◆開発",
  "🔥AI,Ai",
  "◆開発,システム開発",
  "★Java ◇アプリ",
  "◇管理",
  "",
  "",
  "",
  ""
)

$count = $titles.Count
for ($i = 0; $i -lt $count; $i++) {
  $r = $i + 2

  $ws.Cells.Item($r, 1).Value = $timestamp
  $ws.Cells.Item($r, 2).Value = $titles[$i]
  $ws.Cells.Item($r, 3).Value = $categories[$i]
  $ws.Cells.Item($r, 4).Value = $prices[$i]
  $ws.Cells.Item($r, 5).Value = $deadlines[$i]
  $ws.Cells.Item($r, 6).Value = $urls[$i]
  $ws.Cells.Item($r, 7).Value = $scores[$i]

  if ($skills[$i] -ne "") {
    $ws.Cells.Item($r, 8).Value = $skills[$i]
  } else {
    $ws.Cells.Item($r, 8).ClearContents()
  }

  $ws.Hyperlinks.Add($ws.Cells.Item($r, 6), $urls[$i])
  $ws.Cells.Item($r, 6).Style = "Hyperlink"
}
